# Generate Report for Handback
# The handback transform failed for the f5cf148b-... file in both the
# zh-cn and de-de locales. Update the status shown on the Overview sheet
# as well as on each locale-specific sheet, and record the detailed
# error message in the "Error Detail" column (P) for that file's row.

$wb = $excel.ActiveWorkbook

$statusText = "Handback transform failed"

$zhErrorDetail = "Handback file name: sbrqdifs.bza is different with handoff file name: f5cf148b-fb95-41d2-9182-15a4abdcef62.d514dae71453899cae3fbae038f45b6bafa9ff08.zh-cn."
$deErrorDetail = "Handback file name: sbrqdifs.bza is different with handoff file name: f5cf148b-fb95-41d2-9182-15a4abdcef62.d514dae71453899cae3fbae038f45b6bafa9ff08.de-de."

# --- Overview sheet: update the per-locale status for the f5cf148b.md row (row 3)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# --- zh-cn sheet: update Status (C3) and Error Detail (P3) for the f5cf148b.md row
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $statusText
$wsZh.Range("P3").Value = $zhErrorDetail
$wsZh.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: update Status (C3) and Error Detail (P3) for the f5cf148b.md row
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $statusText
$wsDe.Range("P3").Value = $deErrorDetail
$wsDe.Columns.Item(16).ColumnWidth = 39.17
